# This script applies the commit's change: two pairs of data rows in the
# "Artfynd" sheet had their record content swapped (rows 8<->9 and rows
# 17<->18). We read the old values into variables first, then write the
# swapped values back cell-by-cell so we never clobber a value before it
# has been captured.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 8 <-> 9 : only columns A, P, Q, R actually differ between the two
# records; every other cell in the two rows is already identical.
# ---------------------------------------------------------------------

$A8 = $ws.Range("A8").Value2
$P8 = $ws.Range("P8").Value2
$Q8 = $ws.Range("Q8").Value2
$R8 = $ws.Range("R8").Value2

$A9 = $ws.Range("A9").Value2
$P9 = $ws.Range("P9").Value2
$Q9 = $ws.Range("Q9").Value2
$R9 = $ws.Range("R9").Value2

$ws.Range("A8").Value2 = $A9
$ws.Range("P8").Value2 = $P9
$ws.Range("Q8").Value2 = $Q9
$ws.Range("R8").Value2 = $R9

$ws.Range("A9").Value2 = $A8
$ws.Range("P9").Value2 = $P8
$ws.Range("Q9").Value2 = $Q8
$ws.Range("R9").Value2 = $R8

# ---------------------------------------------------------------------
# Rows 17 <-> 18 : the two observation records are fully swapped. Only
# the columns that actually differ are touched below (columns such as
# C, N, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AF, AG, AT, AW, AX, AY hold
# identical values in both rows already, so they are left untouched).
# ---------------------------------------------------------------------

$A17 = $ws.Range("A17").Value2
$B17 = $ws.Range("B17").Value2
$D17 = $ws.Range("D17").Value2
$E17 = $ws.Range("E17").Value2
$F17 = $ws.Range("F17").Value2
$G17 = $ws.Range("G17").Value2
$H17 = $ws.Range("H17").Value2
$I17 = $ws.Range("I17").Value2
$J17 = $ws.Range("J17").Value2
$K17 = $ws.Range("K17").Value2
$P17 = $ws.Range("P17").Value2
$Q17 = $ws.Range("Q17").Value2
$R17 = $ws.Range("R17").Value2
$AJ17 = $ws.Range("AJ17").Value2
$AK17 = $ws.Range("AK17").Value2
$AM17 = $ws.Range("AM17").Value2
$AO17 = $ws.Range("AO17").Value2

$A18 = $ws.Range("A18").Value2
$B18 = $ws.Range("B18").Value2
$D18 = $ws.Range("D18").Value2
$E18 = $ws.Range("E18").Value2
$F18 = $ws.Range("F18").Value2
$G18 = $ws.Range("G18").Value2
$H18 = $ws.Range("H18").Value2
$I18 = $ws.Range("I18").Value2
$J18 = $ws.Range("J18").Value2
$K18 = $ws.Range("K18").Value2
$P18 = $ws.Range("P18").Value2
$Q18 = $ws.Range("Q18").Value2
$R18 = $ws.Range("R18").Value2

# Row 17 becomes the old row 18 content (Knärot / Goodyera repens, VU)
$ws.Range("A17").Value2 = $A18
$ws.Range("B17").Value2 = $B18
$ws.Range("D17").Value2 = $D18
$ws.Range("E17").Value2 = $E18
$ws.Range("F17").Value2 = $F18
$ws.Range("G17").Value2 = $G18
$ws.Range("H17").Value2 = $H18
$ws.Range("I17").Value2 = $I18
$ws.Range("J17").Value2 = $J18
$ws.Range("K17").Value2 = $K18
$ws.Range("P17").Value2 = $P18
$ws.Range("Q17").Value2 = $Q18
$ws.Range("R17").Value2 = $R18
# row 18 had no substrate info, so row 17 loses its AJ/AK/AM/AO values
$ws.Range("AJ17").Value2 = ""
$ws.Range("AK17").Value2 = ""
$ws.Range("AM17").Value2 = ""
$ws.Range("AO17").Value2 = ""

# Row 18 becomes the old row 17 content (Brunpudrad nållav / Chaenotheca
# gracillima, NT) and gains the substrate info that used to live on row 17
$ws.Range("A18").Value2 = $A17
$ws.Range("B18").Value2 = $B17
$ws.Range("D18").Value2 = $D17
$ws.Range("E18").Value2 = $E17
$ws.Range("F18").Value2 = $F17
$ws.Range("G18").Value2 = $G17
$ws.Range("H18").Value2 = $H17
$ws.Range("I18").Value2 = $I17
$ws.Range("J18").Value2 = $J17
$ws.Range("K18").Value2 = $K17
$ws.Range("P18").Value2 = $P17
$ws.Range("Q18").Value2 = $Q17
$ws.Range("R18").Value2 = $R17
$ws.Range("AJ18").Value2 = $AJ17
$ws.Range("AK18").Value2 = $AK17
$ws.Range("AM18").Value2 = $AM17
$ws.Range("AO18").Value2 = $AO17
